$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A (event codes) is treated as text, matching the source data
$ws.Range("A2:A37").NumberFormat = "@"

$ws.Range("A2").Value = '100'
$ws.Range("B2").Value = 'Accidente ofidico'
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = '113'
$ws.Range("B3").Value = 'Desnutrici”n aguda en menores de 5 anos'
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 0.18

$ws.Range("A4").Value = '115'
$ws.Range("B4").Value = 'Cancer en menores de 18 anos'
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.37

$ws.Range("A5").Value = '155'
$ws.Range("B5").Value = 'Cancer de la mama y cuello uterino'
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.09

$ws.Range("A6").Value = '210'
$ws.Range("B6").Value = 'Dengue'
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 0

$ws.Range("A7").Value = '215'
$ws.Range("B7").Value = 'Defectos congenitos'
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.06

$ws.Range("A8").Value = '220'
$ws.Range("B8").Value = 'Dengue grave'
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1

$ws.Range("A9").Value = '300'
$ws.Range("B9").Value = 'Agresiones por animales potencialmente transmisores de rabia'
$ws.Range("C9").Value = 44
$ws.Range("D9").Value = 56
$ws.Range("E9").Value = 0.01

$ws.Range("A10").Value = '330'
$ws.Range("B10").Value = 'Hepatitis a'
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.37

$ws.Range("A11").Value = '340'
$ws.Range("B11").Value = 'Hepatitis b, c y coinfeccion hepatitis b y delta'
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0.37

$ws.Range("A12").Value = '342'
$ws.Range("B12").Value = 'Enfermedades huerfanas - raras'
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 0.16

$ws.Range("A13").Value = '346'
$ws.Range("B13").Value = 'Ira por virus nuevo'
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.14

$ws.Range("A14").Value = '348'
$ws.Range("B14").Value = 'Infeccion respiratoria aguda grave irag inusitada'
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1

$ws.Range("A15").Value = '352'
$ws.Range("B15").Value = 'Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico'
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.37

$ws.Range("A16").Value = '355'
$ws.Range("B16").Value = 'Enfermedad transmitida por alimentos o agua (eta)'
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0

$ws.Range("A17").Value = '356'
$ws.Range("B17").Value = 'Intento de suicidio'
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 0.01

$ws.Range("A18").Value = '357'
$ws.Range("B18").Value = 'Iad - infecciones asociadas a dispositivos - individual'
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0.37

$ws.Range("A19").Value = '365'
$ws.Range("B19").Value = 'Intoxicaciones'
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0.12

$ws.Range("A20").Value = '420'
$ws.Range("B20").Value = 'Leishmaniasis cutanea'
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 1

$ws.Range("A21").Value = '450'
$ws.Range("B21").Value = 'Lepra'
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0

$ws.Range("A22").Value = '455'
$ws.Range("B22").Value = 'Leptospirosis'
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0.27

$ws.Range("A23").Value = '465'
$ws.Range("B23").Value = 'Malaria'
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 1

$ws.Range("A24").Value = '535'
$ws.Range("B24").Value = 'Meningitis bacteriana y enfermedad meningoc”cica'
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 0

$ws.Range("A25").Value = '549'
$ws.Range("B25").Value = 'Morbilidad materna extrema'
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0.18

$ws.Range("A26").Value = '560'
$ws.Range("B26").Value = 'Mortalidad perinatal y neonatal tardia'
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0.27

$ws.Range("A27").Value = '580'
$ws.Range("B27").Value = 'Mortalidad por dengue'
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 1

$ws.Range("A28").Value = '591'
$ws.Range("B28").Value = 'Vigilancia integrada de muertes en menores de cinco anos por infeccion respiratoria aguda - enfermedad diarreica aguda y/o desnutricion'
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 1

$ws.Range("A29").Value = '610'
$ws.Range("B29").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("D29").Value = 0
$ws.Range("E29").ClearContents()

$ws.Range("A30").Value = '620'
$ws.Range("B30").Value = 'Parotiditis'
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0.37

$ws.Range("A31").Value = '740'
$ws.Range("B31").Value = 'Sifilis congenita'
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 1

$ws.Range("A32").Value = '750'
$ws.Range("B32").Value = 'Sifilis gestacional'
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0.18

$ws.Range("A33").Value = '760'
$ws.Range("B33").Value = 'Tetanos accidental'
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 1

$ws.Range("A34").Value = '800'
$ws.Range("B34").Value = 'Tos ferina'
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 1

$ws.Range("A35").Value = '813'
$ws.Range("B35").Value = 'Tuberculosis'
$ws.Range("C35").Value = 7
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 0.15

$ws.Range("A36").Value = '831'
$ws.Range("B36").Value = 'Varicela individual'
$ws.Range("C36").Value = 8
$ws.Range("D36").Value = 5
$ws.Range("E36").Value = 0.09

$ws.Range("A37").Value = '850'
$ws.Range("B37").Value = 'Vih/sida/mortalidad por sida'
$ws.Range("C37").Value = 11
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = 0
